$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $r = $d.Content
    $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "0.97 (0.83 to 1.15), p = 0.772" "0.974 (0.825 to 1.15), p = 0.772"
Replace-Text "0.99 (0.98 to 1), p = 0.187" "0.99 (0.977 to 1.004), p = 0.187"
Replace-Text "1.19 (1.14 to 1.25), p < 0.001" "1.192 (1.138 to 1.248), p < 0.001"
Replace-Text "1.01 (1.01 to 1.01), p < 0.001" "1.009 (1.006 to 1.012), p < 0.001"
Replace-Text "1.28 (1.22 to 1.34), p < 0.001" "1.279 (1.22 to 1.34), p < 0.001"
Replace-Text "0.99 (0.99 to 0.99), p < 0.001" "0.99 (0.987 to 0.993), p < 0.001"
Replace-Text "0.49 (0.17 to 1.39), p = 0.237" "1.206 (0.432 to 3.366), p = 0.744"
Replace-Text "1.02 (0.94 to 1.11), p = 0.635" "0.894 (0.79 to 1.013), p = 0.114"
Replace-Text "1.1 (1.04 to 1.16), p = 0.003" "1.151 (1.09 to 1.215), p < 0.001"
Replace-Text "0.99 (0.98 to 1), p = 0.111" "0.987 (0.983 to 0.991), p < 0.001"
Replace-Text "2.48 (2.37 to 2.6), p < 0.001" "2.483 (2.374 to 2.597), p < 0.001"
Replace-Text "0.99 (0.98 to 0.99), p < 0.001" "0.988 (0.983 to 0.992), p < 0.001"
Replace-Text "1.46 (1.23 to 1.73), p < 0.001" "1.457 (1.23 to 1.726), p < 0.001"
Replace-Text "1 (0.99 to 1.02), p = 0.442" "1.005 (0.993 to 1.017), p = 0.442"
Replace-Text "0.84 (0.73 to 0.97), p = 0.031" "0.844 (0.733 to 0.971), p = 0.031"
Replace-Text "1 (0.99 to 1.01), p = 0.621" "0.997 (0.987 to 1.008), p = 0.621"

# Widen the third grid column from 3254 to 3622 twips (20 twips = 1 point)
$table = $d.Tables.Item(1)
$table.Columns.Item(3).Width = 3622 / 20.0

Write-Output "Done"
